$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.933.95"
$c.ClearFormats()
$ws.Range("E2").Value = "  +0.63%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.554.95"
$c.ClearFormats()
$ws.Range("E3").Value = "  +1.11%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E4").Value = "  +0.37%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "207.25"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.79%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "21.72"
$c.ClearFormats()
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("E10").Value = "  +1.81%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0858"
$c.ClearFormats()
$ws.Range("E11").Value = "  +0.65%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.773.98"
$c.ClearFormats()
$ws.Range("E12").Value = "  +0.86%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.552.76"
$c.ClearFormats()
$ws.Range("E13").Value = "  +0.67%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.73"
$c.ClearFormats()
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("E15").Value = "  +1.90%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "61.90"
$c.ClearFormats()
$ws.Range("E16").Value = "  +1.26%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "26.898.62"
$c.ClearFormats()
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("E19").Value = "  +0.26%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.24"
$c.ClearFormats()
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("E22").Value = "  +0.27%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.18"
$c.ClearFormats()
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("E24").Value = "  +0.37%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "152.34"
$c.ClearFormats()
$ws.Range("E25").Value = "  -0.47%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.65"
$c.ClearFormats()
$ws.Range("E26").Value = "  +2.80%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "14.89"
$c.ClearFormats()
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("E32").Value = "  +0.20%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.419.27"
$c.ClearFormats()
$ws.Range("E33").Value = "  +4.27%  "
$ws.Range("E34").Value = "  +3.23%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.57"
$c.ClearFormats()
$ws.Range("E35").Value = "  +4.07%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.962"
$c.ClearFormats()
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("E38").Value = "  +1.20%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.523"
$c.ClearFormats()
$ws.Range("E39").Value = "  +0.33%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.807"
$c.ClearFormats()
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("E41").Value = "  +0.45%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.64"
$c.ClearFormats()
$ws.Range("E42").Value = "  -1.15%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.989"
$c.ClearFormats()
$ws.Range("E43").Value = "  -0.54%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.28"
$c.ClearFormats()
$ws.Range("E44").Value = "  +3.88%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "63.72"
$c.ClearFormats()
$ws.Range("E46").Value = "  +0.64%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.689.05"
$c.ClearFormats()
$ws.Range("E47").Value = "  +0.89%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "86.18"
$c.ClearFormats()
$ws.Range("E48").Value = "  +0.41%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0517"
$c.ClearFormats()
$ws.Range("E49").Value = "  +1.41%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0958"
$c.ClearFormats()
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").Value = "0.0₇0967"
$ws.Range("E51").Value = "  -0.67%  "
